$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column C for all existing data rows (2-43)
#    from 45730 to 45731.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 3).Value = 45731
}

# 2. Row 43 gains an explicit row height (matches default but becomes "custom").
$ws.Rows.Item(43).RowHeight = 15

# 3. Append the new record as row 44.
$ws.Cells.Item(44, 1).Value = "A 10579-2025"

$ws.Cells.Item(44, 2).Value = 45721
$ws.Cells.Item(44, 2).NumberFormat = $ws.Cells.Item(43, 2).NumberFormat

$ws.Cells.Item(44, 3).Value = 45731
$ws.Cells.Item(44, 3).NumberFormat = $ws.Cells.Item(43, 3).NumberFormat

$ws.Cells.Item(44, 4).Value = "OKÄNT"
$ws.Cells.Item(44, 5).Value = "OKÄNT"

$ws.Cells.Item(44, 7).Value = 0.8
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 0
$ws.Cells.Item(44, 16).Value = 0
$ws.Cells.Item(44, 17).Value = 0

# R44 keeps the same wrap-text style used by the rest of column R.
$ws.Cells.Item(44, 18).WrapText = $true
